$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new log rows (120 and 121) after the existing data (which ends at row 119)
$ws.Range("A120").Value = 119
$ws.Range("B120").Value = 1
$ws.Range("C120").Value = "2024-06-17 11:11:07"
$ws.Range("D120").Value = 200
$ws.Range("E120").Value = 5

$ws.Range("A121").Value = 120
$ws.Range("B121").Value = 2
$ws.Range("C121").Value = "2024-06-17 11:11:08"
$ws.Range("D121").Value = 200
$ws.Range("E121").Value = 0
